$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values (D=Fecha serial, J=Volumen, K=Precio minimo, L=Precio maximo,
# M=Precio promedio ponderado, P=Precio $/Kg) for rows 19-48.
# Rows 21-31 are untouched by this edit.
$rows = @{
    19 = @{ D = 44391; J = 100;  K = 26000; L = 28000; M = 27000; P = 1080 }
    20 = @{ D = 44679; J = 400;  K = 25000; L = 27000; M = 26000; P = 1040 }
    32 = @{ D = 44426; J = 400;  K = 28000; L = 30000; M = 29000; P = 1160 }
    33 = @{ D = 44461; J = 500;  K = 23000; L = 25000; M = 24000; P = 960  }
    34 = @{ D = 44364; J = 200;  K = 28000; L = 30000; M = 29000; P = 1160 }
    35 = @{ D = 44405; J = 500;  K = 26000; L = 28000; M = 27000; P = 1080 }
    36 = @{ D = 44435; J = 900;  K = 28000; L = 30000; M = 29000; P = 1160 }
    37 = @{ D = 44343; J = 200;  K = 26000; L = 28000; M = 27000; P = 1080 }
    38 = @{ D = 44454; J = 1000; K = 28000; L = 30000; M = 29000; P = 1160 }
    39 = @{ D = 44482; J = 500;  K = 18000; L = 20000; M = 19000; P = 760  }
    40 = @{ D = 44398; J = 500;  K = 26000; L = 28000; M = 27000; P = 1080 }
    41 = @{ D = 44370; J = 400;  K = 27000; L = 28000; M = 27500; P = 1100 }
    42 = @{ D = 44406; J = 600;  K = 26000; L = 28000; M = 27000; P = 1080 }
    43 = @{ D = 44455; J = 800;  K = 28000; L = 30000; M = 29000; P = 1160 }
    44 = @{ D = 44448; J = 400;  K = 28000; L = 30000; M = 29000; P = 1160 }
    45 = @{ D = 44497; J = 500;  K = 13000; L = 15000; M = 14000; P = 560  }
    46 = @{ D = 44357; J = 340;  K = 28000; L = 30000; M = 29000; P = 1160 }
    47 = @{ D = 44419; J = 600;  K = 27000; L = 29000; M = 28000; P = 1120 }
    48 = @{ D = 44489; J = 400;  K = 18000; L = 20000; M = 19000; P = 760  }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("P$r").Value = $vals.P
}

# Row 49 (now duplicate of the shifted-up row 48) is dropped entirely.
$ws.Range("A49:R49").EntireRow.Delete()
